$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dialog act annotations (DAMSLTag in column I, DialogAct in column J)
# following a re-run of SGNN after transcript clean-up.

$ws.Range("I14").Value = "sd"
$ws.Range("J14").Value = "Statement-non-opinion"

$ws.Range("I25").Value = "sv"
$ws.Range("J25").Value = "Statement-opinion"

$ws.Range("I26").Value = "aa"
$ws.Range("J26").Value = "Agree/Accept"

$ws.Range("I31").Value = "sv"
$ws.Range("J31").Value = "Statement-opinion"

$ws.Range("I35").Value = "sv"
$ws.Range("J35").Value = "Statement-opinion"

$ws.Range("I51").Value = "sv"
$ws.Range("J51").Value = "Statement-opinion"

$ws.Range("I57").Value = "sd"
$ws.Range("J57").Value = "Statement-non-opinion"

$ws.Range("I65").Value = "%"
$ws.Range("J65").Value = "Uninterpretable"

$ws.Range("I68").Value = "sd"
$ws.Range("J68").Value = "Statement-non-opinion"
